$d = $word.ActiveDocument

# ===================================================================
# Phase 1: all TEXT content edits for this paragraph.
#
# This engine re-normalizes (merges) adjacent same-formatted runs in a
# paragraph every time that paragraph's text is edited, so every pure
# text change has to happen BEFORE any deliberate run-splitting below -
# otherwise a later text edit would silently undo an earlier split.
# ===================================================================

# "Retrieval of tourist data from " -> "Extracting tourist data from "
$d.Content.Find.Execute("Retrieval of", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Extracting", 2) | Out-Null

# "the geographic information database" -> "a spatial information database based on"
$d.Content.Find.Execute("the geographic information database", $true, $false, $false, $false, $false, `
    $true, 1, $false, "a spatial information database based on", 2) | Out-Null

# "OpenStreetMaps" -> "OpenStreetMap"
$d.Content.Find.Execute("OpenStreetMaps", $true, $false, $false, $false, $false, `
    $true, 1, $false, "OpenStreetMap", 2) | Out-Null

# ===================================================================
# Phase 2: re-split runs and relocate the _GoBack bookmark.
#
# A zero-length bookmark placed at a boundary keeps that boundary from
# being merged away by later edits, even after the bookmark itself is
# deleted again - so we use disposable bookmarks purely to pin the run
# boundaries the final XML needs, then remove them.
# ===================================================================

$text = $d.Content.Text
$base = $text.IndexOf("Extracting")

# Relative offsets from $base, derived from the target run layout:
#   "Extracting" | " tourist data from " | "a" | <<bookmark>> | " " |
#   "spatial" | " information database" | " " | "based on " |
#   "Ope"(i) | "nStreetMap"(i) | ": "(i) | "t"
$offBookmark = 30  # right after "a", before the following " "
$splitOffsets = @(10, 29, 31, 38, 59, 60, 69, 72, 82, 84)

# Move (or create) the hidden _GoBack bookmark to the new location.
$bmRng = $d.Range($base + $offBookmark, $base + $offBookmark)
$d.Bookmarks.Add("_GoBack", $bmRng) | Out-Null

# Force a run boundary at every other split point using a throwaway bookmark.
foreach ($o in $splitOffsets) {
    $p = $base + $o
    $r = $d.Range($p, $p)
    $d.Bookmarks.Add("ZZTMPSPLIT", $r) | Out-Null
    $d.Bookmarks("ZZTMPSPLIT").Delete()
}
